# feat: add 2022-Q1 data
#
# The workbook tracks one sheet per quarter plus a running "总计" (total)
# summary sheet as the last tab. A new quarter (2022-Q1) is being added:
#   1. The existing "总计" sheet is renamed to "2022-Q1" and repopulated
#      with that quarter's fund-holder detail rows (same shape as the
#      other quarterly sheets, e.g. "2021-Q4").
#   2. A copy of it becomes the new "总计" sheet (placed right after it),
#      which is then cleared and rebuilt with the rolling summary table
#      plus the new 2022-Q1 row on top of the previously existing ones.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 0: turn the current last sheet ("总计") into "2022-Q1", then
# duplicate it (while it still carries the original sheet's page setup)
# to become the new "总计" tab. Doing the copy before either sheet's
# contents are rewritten keeps both sheets' formatting/pageSetup intact.
# ---------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsQ1.Name = "2022-Q1"

$wsQ1.Copy($null, $wsQ1)
$wsTotal = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTotal.Name = "总计"

# ---------------------------------------------------------------------
# Step 1: rebuild "2022-Q1" with the new quarter's fund-holder detail.
# ---------------------------------------------------------------------
$wsQ1.Cells.Clear()

$headersQ1 = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headersQ1.Length; $i++) {
    $wsQ1.Cells.Item(1, $i + 2).Value = $headersQ1[$i]
}

# Numeric-looking text columns (fund code, size, position, rank %, NAV)
# must stay text even though they look numeric, matching how every other
# quarter sheet stores these values. The fund name (C) is plain text
# already, so it is left on the default General format.
$wsQ1.Range("B2").NumberFormat = "@"
$wsQ1.Range("D2:G2").NumberFormat = "@"

$wsQ1.Cells.Item(2, 1).Value = 0
$wsQ1.Cells.Item(2, 2).Value = "159962"
$wsQ1.Cells.Item(2, 3).Value = "华夏中证四川国企改革ETF"
$wsQ1.Cells.Item(2, 4).Value = "0.49"
$wsQ1.Cells.Item(2, 5).Value = "95.82"
$wsQ1.Cells.Item(2, 6).Value = "3.06"
$wsQ1.Cells.Item(2, 7).Value = "0.0150"
$wsQ1.Cells.Item(2, 8).Value = 7

# Header row + index cell (column A) use the bold/centered/bordered style
# already used by the other quarterly sheets.
$wb.Worksheets.Item(5).Range("B1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122) # xlPasteFormats
$wb.Worksheets.Item(5).Range("A2").Copy()
$wsQ1.Range("A2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Step 2: rebuild "总计" with the updated rolling summary (new 2022-Q1
# row on top, everything else shifted down).
# ---------------------------------------------------------------------
$wsTotal.Cells.Clear()

$wsTotal.Cells.Item(1, 2).Value = "日期"
$wsTotal.Cells.Item(1, 3).Value = "持有数量(只)"
$wsTotal.Cells.Item(1, 4).Value = "持有市值(亿元)"

$totalRows = @(
    @(0, "2022-Q1", 1, 0.02),
    @(1, "2021-Q4", 14, 3.09),
    @(2, "2021-Q3", 30, 15.77),
    @(3, "2021-Q2", 3, 0.03),
    @(4, "2021-Q1", 5, 0.1),
    @(5, "2020-Q4", 1, 0.07000000000000001)
)

foreach ($row in $totalRows) {
    $r = [int]$row[0] + 2
    $wsTotal.Cells.Item($r, 1).Value = $row[0]
    $wsTotal.Cells.Item($r, 2).Value = $row[1]
    $wsTotal.Cells.Item($r, 3).Value = $row[2]
    $wsTotal.Cells.Item($r, 4).Value = $row[3]
}

# Same header/index styling as before.
$wb.Worksheets.Item(5).Range("B1").Copy()
$wsTotal.Range("B1:D1").PasteSpecial(-4122)
$wb.Worksheets.Item(5).Range("A2").Copy()
$wsTotal.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the originally active sheet/selection so we don't introduce an
# unrelated change to the workbook's active-tab bookkeeping.
$wb.Worksheets.Item(1).Activate()

Write-Output "done"
